$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J (copy formatting/style from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I (I0) and J (IF), rows 2-16
$iValues = @(7, 10, 8, 5, 5, 4, 6, 7, 6, 7, 6, 6, 6, 7, 1)
$jValues = @(8, 10, 8, 5, 6, 5, 7, 7, 7, 7, 7, 6, 7, 8, 3)

for ($r = 2; $r -le 16; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
